$d = $word.ActiveDocument

# 1. Replace the font "Times-Roman" with "Times New Roman" throughout the
#    declaration paragraph (ascii + hAnsi). Using Find/Replace with Font
#    formatting so it touches every run carrying that font, regardless of
#    other character formatting (bold/color/etc).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Font.Name = "Times-Roman"
$find.Replacement.Font.Name = "Times New Roman"
$find.Text = ""
$find.Replacement.Text = ""
$find.Execute(
    [ref]"",
    [ref]$false,
    [ref]$false,
    [ref]$false,
    [ref]$false,
    [ref]$false,
    [ref]$true,
    [ref]1,
    [ref]$true,
    [ref]"",
    [ref]2
)

# 2. Word count changed from "sixth" to "seventh".
$d.Content.Find.Execute("sixth", $true, $false, $false, $false, $false,
                         $true, 1, $false, "seventh", 2)

# 3. Merge the split "Vivekananda College of Engineering & " / "Technology,
#    Puttur" run pair into a single run / single piece of text.
$d.Content.Find.Execute("Vivekananda College of Engineering & Technology, Puttur",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false, "Vivekananda College of Engineering & Technology, Puttur", 2)
